$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 6841.3335
$ws.Cells.Item(62, 9).Value = 6668.3335
$ws.Cells.Item(62, 11).Value = 6668.3335
$ws.Cells.Item(62, 13).Value = -6044.3335

$ws.Cells.Item(65, 8).Value = 6841.3335
$ws.Cells.Item(65, 9).Value = 6668.3335
$ws.Cells.Item(65, 11).Value = 33341.6675
$ws.Cells.Item(65, 13).Value = -30221.6675

$ws.Cells.Item(92, 8).Value = 354.91666
$ws.Cells.Item(92, 10).Value = 251.25
$ws.Cells.Item(92, 12).Value = 251.25
$ws.Cells.Item(92, 14).Value = -2747.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4123.7334
$ws.Cells.Item(32, 9).Value = 4123.7334
$ws.Cells.Item(32, 11).Value = 4123.7334
$ws.Cells.Item(32, 13).Value = -3836.7334

$ws.Cells.Item(37, 8).Value = 28500
$ws.Cells.Item(37, 10).Value = 38000
$ws.Cells.Item(37, 12).Value = 38000
$ws.Cells.Item(37, 14).Value = -38546

$ws.Cells.Item(44, 8).Value = 50000
$ws.Cells.Item(44, 10).Value = 50000
$ws.Cells.Item(44, 12).Value = 50000
$ws.Cells.Item(44, 14).Value = -50976

$ws.Cells.Item(51, 8).Value = 49999
$ws.Cells.Item(51, 10).Value = 49999
$ws.Cells.Item(51, 12).Value = 49999
$ws.Cells.Item(51, 14).Value = -51511

$ws.Cells.Item(97, 8).Value = 7636.625
$ws.Cells.Item(97, 9).Value = 156.28572
$ws.Cells.Item(97, 11).Value = 156.28572
$ws.Cells.Item(97, 13).Value = 339.71428

$ws.Cells.Item(132, 8).Value = 1697.7778
$ws.Cells.Item(132, 9).Value = 1658.5
$ws.Cells.Item(132, 11).Value = 4975.5
$ws.Cells.Item(132, 13).Value = -2445.5

$ws.Cells.Item(139, 8).Value = 64895.668
$ws.Cells.Item(139, 10).Value = 64895.668
$ws.Cells.Item(139, 12).Value = 64895.668
$ws.Cells.Item(139, 14).Value = -75175.66800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 38000
$ws.Cells.Item(35, 10).Value = 38000
$ws.Cells.Item(35, 12).Value = 38000
$ws.Cells.Item(35, 14).Value = -38620

$ws.Cells.Item(100, 8).Value = 7962.25
$ws.Cells.Item(100, 10).Value = 7962.25
$ws.Cells.Item(100, 12).Value = 7962.25
$ws.Cells.Item(100, 14).Value = -10126.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(41, 8).Value = 11449.5
$ws.Cells.Item(41, 9).Value = 8085.143
$ws.Cells.Item(41, 11).Value = 8085.143
$ws.Cells.Item(41, 13).Value = -7657.143

$ws.Cells.Item(50, 8).Value = 22500

$ws.Cells.Item(51, 8).Value = 40000
$ws.Cells.Item(51, 10).Value = 40000
$ws.Cells.Item(51, 12).Value = 40000
$ws.Cells.Item(51, 14).Value = -41472

$ws.Cells.Item(59, 8).Value = 44666
$ws.Cells.Item(59, 9).Value = 0
$ws.Cells.Item(59, 10).Value = 44666
$ws.Cells.Item(59, 11).Value = 0
$ws.Cells.Item(59, 12).Value = 44666
$ws.Cells.Item(59, 13).ClearContents()
$ws.Cells.Item(59, 14).Value = -46956

$ws.Cells.Item(60, 8).Value = 23233.166
$ws.Cells.Item(60, 10).Value = 27061.2
$ws.Cells.Item(60, 12).Value = 27061.2
$ws.Cells.Item(60, 14).Value = -28083.2

$ws.Cells.Item(61, 8).Value = 40000
$ws.Cells.Item(61, 10).Value = 40000
$ws.Cells.Item(61, 12).Value = 40000
$ws.Cells.Item(61, 14).Value = -40696

$ws.Cells.Item(68, 8).Value = 47449
$ws.Cells.Item(68, 10).Value = 47449
$ws.Cells.Item(68, 12).Value = 47449
$ws.Cells.Item(68, 14).Value = -48947

$ws.Cells.Item(71, 8).Value = 47449
$ws.Cells.Item(71, 10).Value = 47449
$ws.Cells.Item(71, 12).Value = 142347
$ws.Cells.Item(71, 14).Value = -149835

$ws.Cells.Item(81, 8).Value = 16298
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 14).ClearContents()

$ws.Cells.Item(84, 8).Value = 16298
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 14).ClearContents()

$ws.Cells.Item(86, 8).Value = 13944396
$ws.Cells.Item(86, 9).Value = 17429246
$ws.Cells.Item(86, 11).Value = 17429246
$ws.Cells.Item(86, 13).Value = -17428123

$ws.Cells.Item(88, 8).Value = 17562
$ws.Cells.Item(88, 10).Value = 17562
$ws.Cells.Item(88, 12).Value = 17562
$ws.Cells.Item(88, 14).Value = -18374

$ws.Cells.Item(89, 8).Value = 13944396
$ws.Cells.Item(89, 9).Value = 17429246
$ws.Cells.Item(89, 11).Value = 87146230
$ws.Cells.Item(89, 13).Value = -87140614

$ws.Cells.Item(91, 8).Value = 17562
$ws.Cells.Item(91, 10).Value = 17562
$ws.Cells.Item(91, 12).Value = 17562
$ws.Cells.Item(91, 14).Value = -20370

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(57, 8).Value = 1875
$ws.Cells.Item(57, 9).Value = 500
$ws.Cells.Item(57, 10).Value = 2071.4285
$ws.Cells.Item(57, 11).Value = 1500
$ws.Cells.Item(57, 12).Value = 6214.2855
$ws.Cells.Item(57, 13).Value = -941
$ws.Cells.Item(57, 14).Value = -7332.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(10, 8).Value = 130001.5
$ws.Cells.Item(10, 9).Value = 10003
$ws.Cells.Item(10, 11).Value = 10003
$ws.Cells.Item(10, 13).Value = -9834

$ws.Cells.Item(43, 8).Value = 3374.1667
$ws.Cells.Item(43, 9).Value = 2124.5
$ws.Cells.Item(43, 10).Value = 3999
$ws.Cells.Item(43, 11).Value = 2124.5
$ws.Cells.Item(43, 12).Value = 3999
$ws.Cells.Item(43, 13).Value = -1973.5
$ws.Cells.Item(43, 14).Value = -4301

$ws.Cells.Item(46, 8).Value = 9557
$ws.Cells.Item(46, 9).Value = 3975
$ws.Cells.Item(46, 10).Value = 16999.666
$ws.Cells.Item(46, 11).Value = 3975
$ws.Cells.Item(46, 12).Value = 16999.666
$ws.Cells.Item(46, 13).Value = -3819
$ws.Cells.Item(46, 14).Value = -17311.666

$ws.Cells.Item(57, 9).Value = 20000
$ws.Cells.Item(57, 11).Value = 20000
$ws.Cells.Item(57, 13).Value = -19180

$ws.Cells.Item(122, 8).Value = 2273.0908
$ws.Cells.Item(122, 9).Value = 1889.5555
$ws.Cells.Item(122, 11).Value = 5668.666499999999
$ws.Cells.Item(122, 13).Value = -3218.666499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(25, 8).Value = 13000
$ws.Cells.Item(25, 9).Value = 13000
$ws.Cells.Item(25, 11).Value = 13000
$ws.Cells.Item(25, 13).Value = -12770

$ws.Cells.Item(53, 8).Value = 33681.668
$ws.Cells.Item(53, 9).Value = 25523
$ws.Cells.Item(53, 10).Value = 49999
$ws.Cells.Item(53, 11).Value = 25523
$ws.Cells.Item(53, 12).Value = 49999
$ws.Cells.Item(53, 13).Value = -25005
$ws.Cells.Item(53, 14).Value = -51035

$ws.Cells.Item(56, 8).Value = 7685.3335
$ws.Cells.Item(56, 10).Value = 4057
$ws.Cells.Item(56, 12).Value = 4057
$ws.Cells.Item(56, 14).Value = -5439

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(9, 8).Value = 253
$ws.Cells.Item(9, 9).Value = 253
$ws.Cells.Item(9, 11).Value = 253
$ws.Cells.Item(9, 13).Value = -113

$ws.Cells.Item(31, 8).Value = 19998.5
$ws.Cells.Item(31, 10).Value = 19998.5
$ws.Cells.Item(31, 12).Value = 19998.5
$ws.Cells.Item(31, 14).Value = -20694.5

$ws.Cells.Item(54, 8).Value = 34999.066
$ws.Cells.Item(54, 10).Value = 34999.066
$ws.Cells.Item(54, 12).Value = 34999.066
$ws.Cells.Item(54, 14).Value = -36039.066

$ws.Cells.Item(58, 8).Value = 45093.5
$ws.Cells.Item(58, 9).Value = 45093
$ws.Cells.Item(58, 11).Value = 45093
$ws.Cells.Item(58, 13).Value = -44785

$ws.Cells.Item(92, 8).Value = 58333
$ws.Cells.Item(92, 10).Value = 58333
$ws.Cells.Item(92, 12).Value = 58333
$ws.Cells.Item(92, 14).Value = -63325
